$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("analisis")

# For each data row (2..24) in column C, append " *100" to the existing
# win-rate formula so the value is expressed as a percentage number
# (e.g. 0.8 -> 80) instead of a fraction formatted with a % number format.
for ($row = 2; $row -le 24; $row++) {
    $ws.Range("C$row").Formula = '=COUNTIFS(partidos!C:C, A' + $row + ', partidos!E:E, "Gana")/B' + $row + ' *100'
}

# The cells no longer need the "Porcentaje" (percentage) cell style since
# the formula itself now yields a plain 0-100 number; revert them to the
# default "Normal" style.
$ws.Range("C2:C24").Style = "Normal"

# The "Porcentaje" named cell style is now unused anywhere in the workbook,
# so remove it entirely.
$wb.Styles.Item("Porcentaje").Delete()

# Update the sheet's active selection from D4 to D2.
$ws.Activate()
$ws.Range("D2").Select()
